$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 58
$ws1.Range("F5").Value = 2216
$ws1.Range("F9").Value = 79
$ws1.Range("F10").Value = 68
$ws1.Range("F13").Value = 1901

# Sheet "全部类型" (all types, combined view)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 58
$ws4.Range("F5").Value = 2216
$ws4.Range("F10").Value = 79
$ws4.Range("F11").Value = 68
$ws4.Range("F16").Value = 1901
